$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the student's name (row 7)
$ws.Range("B7").Value = "Cătălina Mădălina Paca"

# Add missing attendance marks (value 2) in column K for several students
$ws.Range("K7").Value = 2
$ws.Range("K10").Value = 2
$ws.Range("K13").Value = 2
$ws.Range("K15").Value = 2
$ws.Range("K20").Value = 2
$ws.Range("K22").Value = 2

# Re-enter the totals formula over the data rows so the attendance sums
# recalculate (this also re-splits the shared formula group, matching how
# Excel regroups Q3:Q22 after the underlying data changed)
$ws.Range("Q3:Q22").Formula = "=SUM(C3:P3)"

# Move the active selection on the frozen pane from M17 to M12
$ws.Range("M12").Select() | Out-Null
